$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 403, pushing all existing
# rows (403..481) down by two (to 405..483).
$ws.Rows.Item(403).Insert()
$ws.Rows.Item(403).Insert()

# New row 403
$ws.Cells.Item(403, 1).Value = 5
$ws.Cells.Item(403, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(403, 3).Value = "Maule"
$ws.Cells.Item(403, 4).Value = 44694
$ws.Cells.Item(403, 5).Value = 7
$ws.Cells.Item(403, 6).Value = "Fruta"
$ws.Cells.Item(403, 7).Value = 100109
$ws.Cells.Item(403, 8).Value = "Uva"
$ws.Cells.Item(403, 9).Value = 100109001
$ws.Cells.Item(403, 10).Value = "Uva"
$ws.Cells.Item(403, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(403, 12).Value = "Especial"
$ws.Cells.Item(403, 13).Value = 100
$ws.Cells.Item(403, 14).Value = 8000
$ws.Cells.Item(403, 15).Value = 8000
$ws.Cells.Item(403, 16).Value = 8000
$ws.Cells.Item(403, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(403, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(403, 19).Value = 444
$ws.Cells.Item(403, 20).Value = 18

# New row 404
$ws.Cells.Item(404, 1).Value = 5
$ws.Cells.Item(404, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(404, 3).Value = "Maule"
$ws.Cells.Item(404, 4).Value = 44694
$ws.Cells.Item(404, 5).Value = 7
$ws.Cells.Item(404, 6).Value = "Fruta"
$ws.Cells.Item(404, 7).Value = 100109
$ws.Cells.Item(404, 8).Value = "Uva"
$ws.Cells.Item(404, 9).Value = 100109001
$ws.Cells.Item(404, 10).Value = "Uva"
$ws.Cells.Item(404, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(404, 12).Value = "Segunda"
$ws.Cells.Item(404, 13).Value = 60
$ws.Cells.Item(404, 14).Value = 6000
$ws.Cells.Item(404, 15).Value = 6000
$ws.Cells.Item(404, 16).Value = 6000
$ws.Cells.Item(404, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(404, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(404, 19).Value = 333
$ws.Cells.Item(404, 20).Value = 18
